# Daily price list refresh for rows 65-83 ("TODAY" sheet)
# Re-applies the updated export data: product rows 65-69 were re-mapped to a
# different product/brand/barcode ordering, and the running SalesQuantity /
# Turnover totals (columns K and L) for the remaining rows were refreshed to
# reflect the new cumulative figures, down to the grand-total row 83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 65: Fitness 375gr -> Fairy Ultra Lemon 400ml ---
$ws.Range("E65").Value = "Fairy® Ultra Lemon 400ml"
$ws.Range("F65").Value = "8001841395883"
$ws.Range("G65").Value = 1.35
$ws.Range("J65").Value = "Fairy"
$ws.Range("L65").Value = 3.04

# --- Row 66: Fairy Ultra Lemon 400ml -> Βερύκοκα Ελληνικά ---
$ws.Range("A66").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("E66").Value = "Βερύκοκα® Ελληνικά (Ζυγιζόμενο) /Kgr"
$ws.Range("F66").Value = "0253"
$ws.Range("G66").Value = 1.95
$ws.Range("H66").Value = 1.65
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = "Βερύκοκα"
$ws.Range("K66").Value = 6
$ws.Range("L66").Value = 7.68

# --- Row 67: Zewa Camomile -> Fitness 375gr ---
$ws.Range("E67").Value = "Fitness® Δημητριακά Απλή Γεύση 375gr"
$ws.Range("F67").Value = "7613034152381"
$ws.Range("G67").Value = 3.65
$ws.Range("J67").Value = "Fitness"
$ws.Range("K67").Value = 6
$ws.Range("L67").Value = 14.53

# --- Row 68: Βερύκοκα Ελληνικά -> Ηλιος Σπαγγέτι Νο10 ---
$ws.Range("A68").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("E68").Value = "Ηλιος® Σπαγγετι Νο10  500gr (2+1Δωρο)"
$ws.Range("F68").Value = "5201020791540"
$ws.Range("G68").Value = 1.78
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 20
$ws.Range("J68").Value = "Ηλιος"
$ws.Range("L68").Value = 7.12

# --- Row 69: Ηλιος Σπαγγέτι Νο10 -> Zewa Camomile ---
$ws.Range("E69").Value = "Zewa® Ρολό Υγείας Camomile Deluxe 8τεμ."
$ws.Range("F69").Value = "7322540055337"
$ws.Range("G69").Value = 4.95
$ws.Range("I69").Value = 30
$ws.Range("J69").Value = "Zewa"
$ws.Range("L69").Value = 16.74

# --- Rows 73-82: refreshed cumulative SalesQuantity / Turnover totals ---
$ws.Range("K73").Value = 17
$ws.Range("L73").Value = 20.56

$ws.Range("K74").Value = 19
$ws.Range("L74").Value = 48.68

$ws.Range("K77").Value = 31.045
$ws.Range("L77").Value = 45.35

$ws.Range("K78").Value = 38
$ws.Range("L78").Value = 33.41

$ws.Range("K79").Value = 42.23
$ws.Range("L79").Value = 31.03

$ws.Range("K80").Value = 51
$ws.Range("L80").Value = 44.89

$ws.Range("K81").Value = 87.28
$ws.Range("L81").Value = 99.33

$ws.Range("H82").Value = 0.85
$ws.Range("K82").Value = 95.435
$ws.Range("L82").Value = 99.73

# --- Row 83: grand totals ---
$ws.Range("K83").Value = 556.99
$ws.Range("L83").Value = 1022.07
